$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The daily-report series had a gap (2021-02-08 / serial 44235 was missing
# between row 92's 44234 and the old row 93's 44236). This update fills in
# that missing day and extends the series by one more day at the end
# (2021-03-02 / serial 44257), which is why the sheet grows from
# A1:D113 to A1:D115.

# Insert a new row at 93, shifting the old rows 93-113 down to 94-114
# (this keeps every shifted cell's formatting, notably column A's date style).
$ws.Rows.Item(93).Insert()

# New row 93 is the filled-in gap day (2021-02-08). Give column A the same
# date style used by the rest of the date column, then fill in the values.
$ws.Cells.Item(92, 1).Copy()
$ws.Cells.Item(93, 1).PasteSpecial(-4122)
$ws.Cells.Item(93, 1).Value = 44235
$ws.Cells.Item(93, 2).Value = 0
$ws.Cells.Item(93, 3).Value = 1
$ws.Cells.Item(93, 4).Value = 46.70714619336758

# Row 92's 7-day rolling totals are recalculated now that the gap day sits
# right after it.
$ws.Cells.Item(92, 3).Value = 0
$ws.Cells.Item(92, 4).Value = 0

# Row 112 (the old row 111, now shifted down by the insert) gets explicit
# rolling totals instead of being blank.
$ws.Cells.Item(112, 3).Value = 0
$ws.Cells.Item(112, 4).Value = 0

# Append a brand-new trailing row 115 for 2021-03-02, matching the style of
# the row above it; its rolling-sum columns stay blank, same as rows 113-114.
$ws.Cells.Item(114, 1).Copy()
$ws.Cells.Item(115, 1).PasteSpecial(-4122)
$ws.Cells.Item(115, 1).Value = 44257
$ws.Cells.Item(115, 2).Value = 0
